$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A52").Value = "dairy milk chocolate 20 rs"
$ws.Range("C52").Value = 20
$ws.Range("D52").Value = "dairy milk chocolate 20 rs.jpg"

$ws.Range("A53").Value = "Sprit 20 rs"
$ws.Range("C53").Value = 10
$ws.Range("D53").Value = "Sprit 20 rs.jpg"

$ws.Range("A54").Value = "5 Star 5 rs"
$ws.Range("A55").Value = "Dite Coke 25 rs"
$ws.Range("C54").Value = 5
$ws.Range("C55").Value = 25
$ws.Range("D54").Value = "5 Star 5 rs.jpg"
$ws.Range("D55").Value = "Dite Coke 25 rs.jpg"

$ws.Range("C52:C55").NumberFormat = $ws.Range("C51").NumberFormat

$ws.Range("D56").Select()
